# Add a new localisation entry (SCRIPT/P01P04A/us2205.ssb) below the existing
# last entry (SCRIPT/P01P04A/us2105.ssb), and append the next filename marker
# row (SCRIPT/P01P04A/us2305.ssb), following the established repeating
# 3-row-group pattern used throughout the sheet:
#   - a "filename" row (bordered/shaded style, columns A:E)
#   - a "data" row (line number + three language variants)
#   - (historically the next filename row starts the following group)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture format templates from an existing, already-styled group -------
# Row 45/46 are the most recent fully formed "filename row" + "data row"
# pair, so we reuse their formatting (borders, fill/font, wrap, row height)
# for the newly inserted rows 47/48. Row 44 gives us the plain single-cell
# "filename only" formatting needed for the final trailing row (49).

# 1) Row 47 currently holds the *old* trailing filename cell (A47, style 4).
#    It becomes a full "filename row" (style 6/7 across A:E), matching the
#    look of row 45.
$ws.Range("A45:E45").Copy()
$ws.Range("A47:E47").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Rows(47).RowHeight = $ws.Rows(45).RowHeight()

# 2) Row 48 becomes the new "data row" (style 4/5 across A:E), matching the
#    look of row 46.
$ws.Range("A46:E46").Copy()
$ws.Range("A48:E48").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Rows(48).RowHeight = $ws.Rows(46).RowHeight()

# 3) Row 49 becomes a lone "filename" cell in column A only (style 4), like
#    row 47 used to be before this edit.
$ws.Range("A44").Copy()
$ws.Range("A49").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Rows(49).RowHeight = $ws.Rows(44).RowHeight()

# --- values ------------------------------------------------------------
# Row 47: unchanged filename text, now styled as a "filename row".
$ws.Range("A47").Value = 'SCRIPT/P01P04A/us2105.ssb'

# Row 48: new localisation entry data.
$ws.Range("A48").Value = 'SCRIPT/P01P04A/us2205.ssb'
$ws.Range("B48").Value = 21
$ws.Range("C48").Value = ' Whooooa!\nI''m ready to explore the whoooole world\nwith [CS:N]Zigzagoon[CR]!'
$ws.Range("D48").Value = ' Вааааау! Я готов исследовать\nвееееесь мир вместе с [CS:N]Зигзагуном[CR]!'
$ws.Range("E48").Value = ' Âàààààô! Ÿ ãïóïâ éòòìåäïâàóû\nâåååååòû íéñ âíåòóå ò [CS:N]Èéãèàãôîïí[CR]!'

# Row 49: next filename marker, trailing row of the sheet.
$ws.Range("A49").Value = 'SCRIPT/P01P04A/us2305.ssb'

# --- final selection (mirrors what Excel leaves behind after data entry) --
[void]$ws.Range("D48").Select()
